$d = $word.ActiveDocument

$replacements = @(
    @("982×3=", "155×9="),
    @("133×6=", "875×6="),
    @("703×5=", "191×4="),
    @("466×9=", "271×5="),
    @("793×4=", "513×9="),
    @("775×7=", "569×3="),
    @("202×4=", "627×6="),
    @("589×4=", "564×5="),
    @("791×7=", "224×8="),
    @("130×4=", "658×4="),
    @("630×7=", "892×4="),
    @("483×3=", "914×5="),
    @("480×4=", "585×9="),
    @("571×6=", "722×6="),
    @("656×8=", "347×9="),
    @("268×6=", "441×9="),
    @("624×8=", "275×3="),
    @("509×9=", "915×7="),
    @("741×2=", "737×5="),
    @("677×4=", "684×5="),
    @("617×3=", "919×9="),
    @("267×3=", "691×6="),
    @("185×7=", "507×5="),
    @("970×6=", "284×5="),
    @("985×6=", "175×4=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2) | Out-Null
}
